$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 69, shifting existing rows 69:164 down to 70:165
$ws.Rows("69:69").Insert()

# Populate the newly inserted row 69 with the new price-report record
$ws.Range("A69").Value = 3
$ws.Range("B69").Value = "Femacal de La Calera"
$ws.Range("C69").Value = "Coquimbo"
$ws.Range("D69").Value = [datetime]"2022-03-18"
$ws.Range("E69").Value = 5
$ws.Range("F69").Value = 100112030
$ws.Range("G69").Value = "Poroto granado"
$ws.Range("H69").Value = "Sin especificar"
$ws.Range("I69").Value = "Primera"
$ws.Range("J69").Value = 73
$ws.Range("K69").Value = 23000
$ws.Range("L69").Value = 24000
$ws.Range("M69").Value = 23521
$ws.Range("N69").Value = "$/malla 25 kilos"
$ws.Range("O69").Value = "Provincia de Petorca"
$ws.Range("P69").Value = 941
$ws.Range("Q69").Value = 25
$ws.Range("R69").Value = "Hortaliza"
